# Sistemazione master test plan
# Fix errors introduced in the "Maintenance test" activity planning:
#  - correct the "Programmazione" phase description for 3 of the tables
#  - drop the stray "Risorse"/"Rischi e assunzioni"/"Ruoli e responsabilità"
#    rows that shouldn't be there

$d = $word.ActiveDocument
$tables = $d.Tables

# --- Table 1 (Stesura component test plan): drop trailing Risorse +
#     Ruoli e responsabilità rows (rows 6 and 7) ---
$t1 = $tables.Item(1)
$t1.Rows.Item(7).Delete()
$t1.Rows.Item(6).Delete()

# --- Table 2 (Individuazione component test): drop trailing Risorse +
#     Ruoli e responsabilità rows (rows 6 and 7) ---
$t2 = $tables.Item(2)
$t2.Rows.Item(7).Delete()
$t2.Rows.Item(6).Delete()

# --- Table 3 (Implementazione component test): fix phase text, then
#     drop trailing Rischi e assunzioni + Ruoli e responsabilità rows ---
# (Wrap=wdFindStop / Replace=wdReplaceOne so the Find stays inside this
#  table's range instead of touching the other tables' identical text)
$t3 = $tables.Item(3)
$t3.Range.Find.Execute("Da ultimare durante la fase di design della manutenzione.", $true, $false, $false, $false, $false, $true, 0, $false, "Da ultimare durante la fase di implementazione della manutenzione.", 1) | Out-Null
$t3.Rows.Item(8).Delete()
$t3.Rows.Item(7).Delete()

# --- Table 4 (Esecuzione component test): fix phase text, then drop
#     trailing Rischi e assunzioni + Ruoli e responsabilità rows ---
$t4 = $tables.Item(4)
$t4.Range.Find.Execute("Da ultimare durante la fase di design della manutenzione.", $true, $false, $false, $false, $false, $true, 0, $false, "Da ultimare durante la fase di system testing della manutenzione.", 1) | Out-Null
$t4.Rows.Item(8).Delete()
$t4.Rows.Item(7).Delete()

# --- Table 5 (Report component test): fix phase text, then drop the
#     trailing Ruoli e responsabilità row ---
$t5 = $tables.Item(5)
$t5.Range.Find.Execute("Da ultimare durante la fase di design della manutenzione.", $true, $false, $false, $false, $false, $true, 0, $false, "Da ultimare durante la fase di system testing della manutenzione.", 1) | Out-Null
$t5.Rows.Item(7).Delete()

# --- Table 6 (Test task iteration): drop the trailing
#     Ruoli e responsabilità row (no phase-text change needed here) ---
$t6 = $tables.Item(6)
$t6.Rows.Item(7).Delete()
